$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of test data: Testcase4 / (blank) / "Please enter a search term."
$ws.Range("A5").Value = "Testcase4"
$ws.Range("C5").Value = "Please enter a search term."

# Match the author's final selection/active cell
$ws.Range("C5").Select()
